$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing D2 value (46 -> 52)
$ws.Range("D2").Value = 52

# New header row (I1:T1) - shared strings referencing the new categories
$ws.Range("I1").Value = "HD_U1"
$ws.Range("J1").Value = "OT_U1"
$ws.Range("K1").Value = "OOP_U1"
$ws.Range("L1").Value = "HD_U2"
$ws.Range("M1").Value = "OT_U2"
$ws.Range("N1").Value = "OOP_U2"
$ws.Range("O1").Value = "HD_U3"
$ws.Range("P1").Value = "OT_U3"
$ws.Range("Q1").Value = "OOP_U3"
$ws.Range("R1").Value = "HD_LPV"
$ws.Range("S1").Value = "OT_LPV"
$ws.Range("T1").Value = "OOP_LPV"

# New data row (I2:T2)
$ws.Range("I2").Value = 200
$ws.Range("J2").Value = 365
$ws.Range("K2").Value = 245
$ws.Range("L2").Value = 352
$ws.Range("M2").Value = 261
$ws.Range("N2").Value = 109
$ws.Range("O2").Value = 100
$ws.Range("P2").Value = 65
$ws.Range("Q2").Value = 156
$ws.Range("R2").Value = 100
$ws.Range("S2").Value = 100
$ws.Range("T2").Value = 606

# Apply header style (s="3") to new header cells, matching existing headers
$ws.Range("H1").Copy()
$ws.Range("I1:T1").PasteSpecial(-4122)

# Apply centered alignment to new data cells
$ws.Range("I2:T2").HorizontalAlignment = -4108

# Update sheet view / selection to match target
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("K3:T3").Select()
